# This script reproduces the commit "Reverted and validated the git conflicts"
# which swaps the old Plant/Nimda test data on the PostDetails (sheet2) and
# CreateSTP_Mandatory (sheet3) sheets for new Swift/Ferrari test data, extends
# the data ranges, and updates the active sheet/selection state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "PostDetails" (2nd sheet): replace Plant1-5 / "Details entered
# successfully for plantN" rows (2-6) with Swift1-20 / "Details entered
# successfully for SwiftN" rows (2-21).
# ---------------------------------------------------------------------------
$wsPost = $wb.Worksheets.Item("PostDetails")

for ($i = 1; $i -le 20; $i++) {
    $row = $i + 1
    $wsPost.Cells.Item($row, 1).Value = "Swift$i"
    $wsPost.Cells.Item($row, 2).Value = "Details entered successfully for Swift$i"
}

# ---------------------------------------------------------------------------
# Sheet "CreateSTP_Mandatory" (3rd sheet): replace Nimda/TEst/Test rows (2-6)
# with Ferrari1-18 / Ferrari123-140 rows (2-19), and extend the
# CommunityOrganizer ("sharath sethu") column down through row 57.
# ---------------------------------------------------------------------------
$wsMand = $wb.Worksheets.Item("CreateSTP_Mandatory")

$mandatoryMsg = 'This is a valid Mandatory Fields scenario. It is supposed give a toast message "STP created successfully!"'
$organizer = "sharath sethu"

for ($i = 1; $i -le 18; $i++) {
    $row = $i + 1
    $wsMand.Cells.Item($row, 1).Value = "Ferrari$i"
    $wsMand.Cells.Item($row, 2).Value = "Ferrari$(122 + $i)"
    $wsMand.Cells.Item($row, 3).Value = $mandatoryMsg
    $wsMand.Cells.Item($row, 4).Value = $organizer
}

# Row 20 keeps the mandatory-message in column C plus the organizer in D
$wsMand.Cells.Item(20, 3).Value = $mandatoryMsg
$wsMand.Cells.Item(20, 4).Value = $organizer

# Rows 21-57 only carry the CommunityOrganizer value in column D
for ($row = 21; $row -le 57; $row++) {
    $wsMand.Cells.Item($row, 4).Value = $organizer
}

# ---------------------------------------------------------------------------
# Selection / active-sheet state.
# Final diff shows CreateSTP_Mandatory's selection moved to D2:D57 and
# PostDetails became the active ("tabSelected") sheet with selection A17:A21.
# Select CreateSTP_Mandatory first so PostDetails ends up as the active tab.
# ---------------------------------------------------------------------------
$wsMand.Activate()
$wsMand.Range("D2:D57").Select()

$wsPost.Activate()
$wsPost.Range("A17:A21").Select()
